# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-24 06:24:17
# Reorders the "Recorded By" (column G) email/name lists on the
# "Session Analysis Results" sheet for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "backup@backdoor.com, system, System"
$ws.Range("G3").Value = "dnasr281@gmail.com, System"
$ws.Range("G5").Value = "backup@backdoor.com, System"
$ws.Range("G6").Value = "dnasr281@gmail.com, System"
$ws.Range("G7").Value = "admin@admin.com, System"
$ws.Range("G8").Value = "backup@backdoor.com, System"

$ws.Range("G28").Value = "backup@backdoor.com, system, System"
$ws.Range("G29").Value = "dnasr281@gmail.com, System"
$ws.Range("G31").Value = "backup@backdoor.com, System"
$ws.Range("G32").Value = "dnasr281@gmail.com, System"
$ws.Range("G33").Value = "admin@admin.com, System"
$ws.Range("G34").Value = "backup@backdoor.com, System"

$ws.Range("G54").Value = "backup@backdoor.com, system, System"
$ws.Range("G55").Value = "dnasr281@gmail.com, System"
$ws.Range("G57").Value = "backup@backdoor.com, System"
$ws.Range("G58").Value = "dnasr281@gmail.com, System"
$ws.Range("G59").Value = "admin@admin.com, System"
$ws.Range("G60").Value = "backup@backdoor.com, System"

$ws.Range("G80").Value = "backup@backdoor.com, System"
$ws.Range("G81").Value = "backup@backdoor.com, System"
$ws.Range("G82").Value = "backup@backdoor.com, System"

$ws.Range("G87").Value = "admin@admin.com, dnasr281@gmail.com"

$ws.Range("G106").Value = "backup@backdoor.com, System"
$ws.Range("G107").Value = "backup@backdoor.com, System"
$ws.Range("G108").Value = "backup@backdoor.com, System"

$ws.Range("G113").Value = "admin@admin.com, dnasr281@gmail.com"

$ws.Range("G132").Value = "backup@backdoor.com, System"
$ws.Range("G133").Value = "backup@backdoor.com, System"
$ws.Range("G134").Value = "backup@backdoor.com, System"

$ws.Range("G139").Value = "admin@admin.com, dnasr281@gmail.com"
